$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.778.78'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '2.102.80'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.62'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.616'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '62.39'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.19%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +2.05%  '
$ws.Range('E10').Value = '  +0.32%  '
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.67'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +5.69%  '
$ws.Range('D13').Value = '2.414.54'
$ws.Range('E13').Value = '  +0.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.05'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.813'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.84%  '
$ws.Range('E16').Value = '  +1.51%  '
$ws.Range('D17').Value = '2.091.13'
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('D18').Value = '38.781.53'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.14'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.66'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('D21').Value = '0.0₃0843'
$ws.Range('E21').Value = '  +1.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.55'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.26%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.34'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.31'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.67'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '171.85'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.137'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.33%  '
$ws.Range('E29').Value = '  +4.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.33'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.53'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +7.32%  '
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('E33').Value = '  +0.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.75'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.07'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +9.93%  '
$ws.Range('E36').Value = '  +1.59%  '
$ws.Range('E37').Value = '  +0.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.55'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.80%  '
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.14'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '102.83'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.90%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0228'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.58%  '
$ws.Range('D43').Value = '1.528.83'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.21'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +7.98%  '
$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.81'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.78%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0918'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.34%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.81'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.08'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +4.46%  '
$ws.Range('E49').Value = '  -0.23%  '
$ws.Range('E50').Value = '  -0.38%  '
$ws.Range('D51').Value = '2.301.06'
$ws.Range('E51').Value = '  +0.59%  '
